$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values (shared strings content change)
$ws.Range("A2").Value = "SWR"
$ws.Range("B2").Value = "Shoftware"

# Update selection to B3
$ws.Range("B3").Select()
